$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column N (GST Number), shifting existing columns right.
$ws.Range("N1").EntireColumn.Insert()

# Populate the new PAN column placeholder value and header (row2 first, matching
# the order the shared strings were authored in).
$ws.Range("N2").Value = "{vendor:pan_no}"
$ws.Range("N1").Value = "PAN"

# Match the formatting of the neighboring header/value cells for the new column.
$ws.Range("N1").Font.Bold = $true
$ws.Range("N1").HorizontalAlignment = -4108

$ws.Range("N2").Font.Bold = $false
$ws.Range("N2").HorizontalAlignment = -4108

# Reflect the cursor landing on the new header cell, as in the authored edit.
[void]$ws.Range("O1").Select()
